$d = $word.ActiveDocument
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $ltype = $p.Range.ListFormat.ListType
    if ($ltype -ne 0) {
        Write-Output "$i : $ltype : $($p.Range.Text)"
    }
}
